$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 updates
$ws.Range("G7").Value = 2.25
$ws.Range("I7").Value = 3.4
$ws.Range("W7").Value = 6.5
$ws.Range("AE7").Value = 17
$ws.Range("AN7").Value = 4

# Row 8 updates
$ws.Range("M8").Value = 1.03
$ws.Range("N8").Value = 11
$ws.Range("U8").Value = 1.9
$ws.Range("V8").Value = 1.86

# Row 9 updates
$ws.Range("M9").Value = 1.17
$ws.Range("N9").Value = 5
$ws.Range("O9").Value = 1.73
$ws.Range("P9").Value = 2
$ws.Range("S9").Value = 1.75
$ws.Range("T9").Value = 2.05
